$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiment2")

# --- Insert two new rows after row 17 (pushes old header row 19.. down to 21..) ---
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# --- Update header labels (comma-space -> underscore style) ---
$ws.Range("B2").Value = "fullProduct_overdraft"
$ws.Range("B21").Value = "fullProduct,_credit"

# --- New summary rows for block 1 (rows 17-19) ---
$ws.Range("B17").Formula = "=AVERAGE(B3:B16)"
$ws.Range("B18").Formula = "=MAX(B3:B16)"
$ws.Range("B19").Formula = "=MIN(B3:B16)"

# Match number formatting/styles of neighboring cells
$ws.Range("B17").NumberFormat = "0.00"
$ws.Range("B18").NumberFormat = "0.00"
$ws.Range("B19").NumberFormat = "0.00"

# --- New summary row for block 2 (row 36, right after row 35) ---
$ws.Range("B36").Formula = "=AVERAGE(B22:B35)"
$ws.Range("B36").NumberFormat = "0.00"

$wb.Application.Calculate()

for ($r = 1; $r -le 36; $r++) {
    $a = $ws.Cells.Item($r, 1).Text
    $b = $ws.Cells.Item($r, 2).Text
    Write-Output "Row $r : A=[$a] B=[$b]"
}
